# Generate Report for Handback
# Applies the "handback" report generation changes to the localization-status
# workbook: updates status text, fills in the "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns for the zh-cn
# and de-de sheets, adds hyperlinks for the newly populated target-file cells,
# and widens a few columns to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: update the per-language status column -----------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# The "Status" column (C) on each language sheet shares the same underlying
# text as the Overview status cells above, so it needs the same update.
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# --- zh-cn sheet: fill target/handback columns for both rows ----------------
$zhTargetFile   = "740ec174-6d4a-4531-ae7d-5e19eca46094.md"
$zhHandbackFile = "740ec174-6d4a-4531-ae7d-5e19eca46094.76a4fffd6dc20842fbee121cd9345c0f6cdeb40b.zh-cn.xlf"
$zhHandbackTime = "2016-11-09 00:56:41"

$zhcn.Range("I2").Value = $zhTargetFile
$zhcn.Range("J2").Value = $zhHandbackFile
$zhcn.Range("K2").Value = $zhHandbackTime

$zhcn.Range("I3").Value = $zhTargetFile
$zhcn.Range("J3").Value = $zhHandbackFile
$zhcn.Range("K3").Value = $zhHandbackTime

# --- de-de sheet: fill target/handback columns for both rows ----------------
$deTargetFile   = "740ec174-6d4a-4531-ae7d-5e19eca46094.md"
$deHandbackFile = "740ec174-6d4a-4531-ae7d-5e19eca46094.76a4fffd6dc20842fbee121cd9345c0f6cdeb40b.de-de.xlf"
$deHandbackTime = "2016-11-09 00:57:00"

$dede.Range("I2").Value = $deTargetFile
$dede.Range("J2").Value = $deHandbackFile
$dede.Range("K2").Value = $deHandbackTime

$dede.Range("I3").Value = $deTargetFile
$dede.Range("J3").Value = $deHandbackFile
$dede.Range("K3").Value = $deHandbackTime

# --- Hyperlinks --------------------------------------------------------------
# Recreate the hyperlinks on both worksheets so that the newly-populated
# "Latest Target File" cells (column I) get a link to the source markdown
# file, in the same row-major order the rows were generated.
$mdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a87e154730d83509be1003f5d227855f346f614/e2e/740ec174-6d4a-4531-ae7d-5e19eca46094.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a87e154730d83509be1003f5d227855f346f614/e2e/ffffe684983d-7fd2-423b-a493-8c85ed6d4b24.md"

foreach ($ws in @($zhcn, $dede)) {
    $ws.Hyperlinks.Delete() | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, [Type]::Missing, [Type]::Missing, "740ec174-6d4a-4531-ae7d-5e19eca46094.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, "740ec174-6d4a-4531-ae7d-5e19eca46094.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, [Type]::Missing, [Type]::Missing, "ffffe684983d-7fd2-423b-a493-8c85ed6d4b24.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrl, [Type]::Missing, [Type]::Missing, "740ec174-6d4a-4531-ae7d-5e19eca46094.md") | Out-Null
}

# --- Column widths ------------------------------------------------------------
# Widen the status column on the Overview sheet (E, F) and the Status (C) /
# Latest Target File (I) / Latest Handback File (J) columns on the language
# sheets to fit the longer text that was just written into them.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}
